# Apply cryptos list update (Fri Jan 26 18:28:11 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "41.920.23"
$ws.Cells.Item(2, 5).Value = "  +4.89%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.265.90"
$ws.Cells.Item(3, 5).Value = "  +2.08%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.17%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "302.32"
$ws.Cells.Item(5, 5).Value = "  +3.41%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  +6.33%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.533"
$ws.Cells.Item(7, 5).Value = "  +4.07%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.999"
$ws.Cells.Item(8, 5).Value = "  -0.16%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.487"
$ws.Cells.Item(9, 5).Value = "  +3.82%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "32.87"
$ws.Cells.Item(10, 5).Value = "  +6.18%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "54.74"
$ws.Cells.Item(11, 5).Value = "  +9.55%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0801"
$ws.Cells.Item(12, 5).Value = "  +2.68%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  +3.13%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  +3.33%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "2.614.95"
$ws.Cells.Item(15, 5).Value = "  +2.17%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  +3.02%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "2.299.25"
$ws.Cells.Item(17, 5).Value = "  +2.14%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  +3.21%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "41.829.95"
$ws.Cells.Item(19, 5).Value = "  +4.84%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "12.22"
$ws.Cells.Item(20, 5).Value = "  +8.12%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "0.0₃0909"
$ws.Cells.Item(21, 5).Value = "  +2.68%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.96"
$ws.Cells.Item(22, 5).Value = "  +3.25%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "67.41"
$ws.Cells.Item(23, 5).Value = "  +2.67%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "242.30"
$ws.Cells.Item(24, 5).Value = "  +1.95%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.60"
$ws.Cells.Item(25, 5).Value = "  +5.87%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.999"
$ws.Cells.Item(26, 5).Value = "  -0.13%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +4.69%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "23.96"
$ws.Cells.Item(28, 5).Value = "  +2.33%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +1.64%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "9.71"
$ws.Cells.Item(30, 5).Value = "  +4.95%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "34.10"
$ws.Cells.Item(31, 5).Value = "  +6.59%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "158.01"
$ws.Cells.Item(32, 5).Value = "  +0.71%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.998"
$ws.Cells.Item(33, 5).Value = "  -0.12%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +4.73%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.0743"
$ws.Cells.Item(35, 5).Value = "  +4.31%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "3.08"
$ws.Cells.Item(36, 5).Value = "  +5.32%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +3.34%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +6.24%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "Celestia"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "16.59"
$ws.Cells.Item(39, 5).Value = "  +7.74%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "Stellar"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.116"
$ws.Cells.Item(40, 5).Value = "  +3.89%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +4.91%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.98"
$ws.Cells.Item(42, 5).Value = "  +6.31%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "20.26"
$ws.Cells.Item(43, 5).Value = "  +13.86%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "2.050.29"
$ws.Cells.Item(44, 5).Value = "  -3.21%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  +3.74%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "10.09"
$ws.Cells.Item(46, 5).Value = "  +1.25%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.92"
$ws.Cells.Item(47, 5).Value = "  +7.67%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -4.50%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "2.490.50"
$ws.Cells.Item(49, 5).Value = "  +2.59%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  +2.14%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.15"
$ws.Cells.Item(51, 5).Value = "  +4.47%  "

